$wb = $excel.ActiveWorkbook

# --- Layer0 sheet ---
$ws0 = $wb.Worksheets.Item("Layer0")
$ws0.Range("B2").Value = -1.352875235549049
$ws0.Range("C2").Value = -0.8778542889328459
$ws0.Range("B3").Value = 1.145162507094057
$ws0.Range("C3").Value = 0.8371067152591948
$ws0.Range("B4").Value = 1.098998930159426
$ws0.Range("C4").Value = -1.132010222760728

# --- Layer1 sheet ---
$ws1 = $wb.Worksheets.Item("Layer1")
$ws1.Range("B2").Value = -1.66174642026249
$ws1.Range("C2").Value = 0.1115069638636245
$ws1.Range("B3").Value = 1.699690975204034
$ws1.Range("C3").Value = -0.6917408396953876
$ws1.Range("B4").Value = -0.2239801905418986
$ws1.Range("C4").Value = 0.6741884741343568
